$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 111 (ALC)
$ws.Range("H111").Value = 305.07693
$ws.Range("I111").Value = 327.4
$ws.Range("J111").Value = 230.66667
$ws.Range("K111").Value = 982.1999999999999
$ws.Range("L111").Value = 692.00001
$ws.Range("M111").Value = 2084.8
$ws.Range("N111").Value = -6826.00001

# Row 115 (ALC)
$ws.Range("H115").Value = 721.25
$ws.Range("I115").Value = 721.25
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2163.75
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -596.75
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 920.4857
$ws.Range("I2").Value = 912.9
$ws.Range("J2").Value = 930.6
$ws.Range("K2").Value = 912.9
$ws.Range("L2").Value = 930.6
$ws.Range("M2").Value = -799.9
$ws.Range("N2").Value = -1156.6

# Row 32 (ARM)
$ws.Range("H32").Value = 11595.031
$ws.Range("I32").Value = 9900.102000000001
$ws.Range("J32").Value = 33144.855
$ws.Range("K32").Value = 9900.102000000001
$ws.Range("L32").Value = 33144.855
$ws.Range("M32").Value = -9613.102000000001
$ws.Range("N32").Value = -33718.855

# Row 54 (ARM)
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 61 (ARM)
$ws.Range("H61").Value = 324513.1
$ws.Range("I61").Value = 1405.1305
$ws.Range("J61").Value = 1253448.5
$ws.Range("K61").Value = 1405.1305
$ws.Range("L61").Value = 1253448.5
$ws.Range("M61").Value = -1193.1305
$ws.Range("N61").Value = -1253872.5

# Row 74 (ARM)
$ws.Range("H74").Value = 4254.1816
$ws.Range("I74").Value = 1098.0714
$ws.Range("J74").Value = 21928.4
$ws.Range("K74").Value = 1098.0714
$ws.Range("L74").Value = 21928.4
$ws.Range("M74").Value = -224.0714
$ws.Range("N74").Value = -23676.4

# Row 77 (ARM)
$ws.Range("H77").Value = 4254.1816
$ws.Range("I77").Value = 1098.0714
$ws.Range("J77").Value = 21928.4
$ws.Range("K77").Value = 5490.357
$ws.Range("L77").Value = 109642
$ws.Range("M77").Value = -1122.357
$ws.Range("N77").Value = -118378

# Row 116 (ARM)
$ws.Range("H116").Value = 920.4857
$ws.Range("I116").Value = 912.9
$ws.Range("J116").Value = 930.6
$ws.Range("K116").Value = 912.9
$ws.Range("L116").Value = 930.6
$ws.Range("M116").Value = 1381.1
$ws.Range("N116").Value = -5518.6

# Row 132 (ARM)
$ws.Range("H132").Value = 6060.2046
$ws.Range("I132").Value = 4163.727
$ws.Range("K132").Value = 12491.181
$ws.Range("M132").Value = -9961.181

# Row 136 (ARM)
$ws.Range("H136").Value = 324513.1
$ws.Range("I136").Value = 1405.1305
$ws.Range("J136").Value = 1253448.5
$ws.Range("K136").Value = 4215.3915
$ws.Range("L136").Value = 3760345.5
$ws.Range("M136").Value = -1665.3915
$ws.Range("N136").Value = -3765445.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 920.4857
$ws.Range("I3").Value = 912.9
$ws.Range("J3").Value = 930.6
$ws.Range("K3").Value = 912.9
$ws.Range("L3").Value = 930.6
$ws.Range("M3").Value = -798.9
$ws.Range("N3").Value = -1158.6

# Row 107 (BSM)
$ws.Range("H107").Value = 944.175
$ws.Range("I107").Value = 970.4857
$ws.Range("K107").Value = 970.4857
$ws.Range("M107").Value = 949.5143

# Row 134 (BSM)
$ws.Range("H134").Value = 2913.6924
$ws.Range("I134").Value = 1748
$ws.Range("J134").Value = 4273.6665
$ws.Range("K134").Value = 5244
$ws.Range("L134").Value = 12820.9995
$ws.Range("M134").Value = -2709
$ws.Range("N134").Value = -17890.9995

$ws = $wb.Worksheets.Item("CRP")
# Row 45 (CRP)
$ws.Range("H45").Value = 12500
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 12500
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 12500
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -13686

# Row 132 (CRP)
$ws.Range("H132").Value = 6180.75
$ws.Range("I132").Value = 6638.7896
$ws.Range("K132").Value = 19916.3688
$ws.Range("M132").Value = -17386.3688

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 767.55884
$ws.Range("I5").Value = 492.82144
$ws.Range("J5").Value = 2049.6667
$ws.Range("K5").Value = 1478.46432
$ws.Range("L5").Value = 6149.000100000001
$ws.Range("M5").Value = -1366.46432
$ws.Range("N5").Value = -6373.000100000001

# Row 39 (CUL)
$ws.Range("H39").Value = 1983.4445
$ws.Range("I39").Value = 650
$ws.Range("J39").Value = 2150.125
$ws.Range("K39").Value = 1950
$ws.Range("L39").Value = 6450.375
$ws.Range("M39").Value = -1656
$ws.Range("N39").Value = -7038.375

# Row 130 (CUL)
$ws.Range("H130").Value = 1723.6
$ws.Range("I130").Value = 872.6667
$ws.Range("J130").Value = 3000
$ws.Range("K130").Value = 2618.0001
$ws.Range("L130").Value = 9000
$ws.Range("M130").Value = 2401.9999
$ws.Range("N130").Value = -19040

# Row 135 (CUL)
$ws.Range("H135").Value = 767.55884
$ws.Range("I135").Value = 492.82144
$ws.Range("J135").Value = 2049.6667
$ws.Range("K135").Value = 4435.39296
$ws.Range("L135").Value = 18447.0003
$ws.Range("M135").Value = -1900.39296
$ws.Range("N135").Value = -23517.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (GSM)
$ws.Range("H113").Value = 1001181.1
$ws.Range("I113").Value = 1667701.9
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1667701.9
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = -1665531.9
$ws.Range("N113").Value = -5740

# Row 132 (GSM)
$ws.Range("H132").Value = 5526.6855
$ws.Range("I132").Value = 5730.5713
$ws.Range("J132").Value = 4711.143
$ws.Range("K132").Value = 17191.7139
$ws.Range("L132").Value = 14133.429
$ws.Range("M132").Value = -14661.7139
$ws.Range("N132").Value = -19193.429

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 1746.4615
$ws.Range("I7").Value = 1300.8
$ws.Range("J7").Value = 2025
$ws.Range("K7").Value = 1300.8
$ws.Range("L7").Value = 2025
$ws.Range("M7").Value = -1188.8
$ws.Range("N7").Value = -2249

# Row 104 (LTW)
$ws.Range("H104").Value = 12552.728
$ws.Range("J104").Value = 12552.728
$ws.Range("L104").Value = 12552.728
$ws.Range("N104").Value = -19540.728

# Row 126 (LTW)
$ws.Range("H126").Value = 1746.4615
$ws.Range("I126").Value = 1300.8
$ws.Range("J126").Value = 2025
$ws.Range("K126").Value = 3902.4
$ws.Range("L126").Value = 6075
$ws.Range("M126").Value = -1432.4
$ws.Range("N126").Value = -11015

$ws = $wb.Worksheets.Item("WVR")
# Row 104 (WVR)
$ws.Range("H104").Value = 34570
$ws.Range("J104").Value = 34570
$ws.Range("L104").Value = 34570
$ws.Range("N104").Value = -41558

# Row 107 (WVR)
$ws.Range("H107").Value = 233.21739
$ws.Range("I107").Value = 180.58824
$ws.Range("J107").Value = 382.33334
$ws.Range("K107").Value = 541.76472
$ws.Range("L107").Value = 1147.00002
$ws.Range("M107").Value = 1378.23528
$ws.Range("N107").Value = -4987.000019999999

# Row 109 (WVR)
$ws.Range("H109").Value = 31300
$ws.Range("J109").Value = 31300
$ws.Range("L109").Value = 31300
$ws.Range("N109").Value = -34074
